# Swap the entire contents of row 40 and row 41 (A:CE) — the Monday and
# Tuesday lab results got entered in the wrong order, so swap them back.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = 83   # column CE
$row1 = 40
$row2 = 41

$vals1 = @()
$vals2 = @()

for ($col = 1; $col -le $lastCol; $col++) {
    $c1 = $ws.Cells.Item($row1, $col)
    $c2 = $ws.Cells.Item($row2, $col)
    $vals1 += ,$c1.Value2
    $vals2 += ,$c2.Value2
}

for ($col = 1; $col -le $lastCol; $col++) {
    $old1 = $vals1[$col - 1]
    $old2 = $vals2[$col - 1]
    # Skip no-op writes (also sidesteps clearing cells that hold an
    # empty-string shared value, since Value = "" blanks a cell in Excel).
    if ($old1 -ne $old2) {
        $ws.Cells.Item($row1, $col).Value = $old2
        $ws.Cells.Item($row2, $col).Value = $old1
    }
}
